# The workbook records daily fruit/vegetable prices, one row per day,
# ordered with the most recent day first (row 2 is the newest record).
# A new daily record (for a date after the current newest one) is being
# added. This pushes the existing rows 18-132 down by one row (to 19-133)
# and inserts the new record's data into row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18:132 down to 19:133, making room for the new record at row 18.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new daily record.
$ws.Cells.Item(18, 1).Value = 11
$ws.Cells.Item(18, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(18, 3).Value = "Bíobío"
$ws.Cells.Item(18, 4).Value = 44831
$ws.Cells.Item(18, 5).Value = 8
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100108
$ws.Cells.Item(18, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(18, 9).Value = 100108002
$ws.Cells.Item(18, 10).Value = "Mango"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 100
$ws.Cells.Item(18, 14).Value = 8000
$ws.Cells.Item(18, 15).Value = 8500
$ws.Cells.Item(18, 16).Value = 8250
$ws.Cells.Item(18, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(18, 18).Value = "Brasil"
$ws.Cells.Item(18, 19).Value = 2062
$ws.Cells.Item(18, 20).Value = 4
